$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''28.449.79'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '''1.564.47'
$ws.Range('E3').Value = '  -1.23%  '
$ws.Range('E4').Value = '  +0.57%  '
$ws.Range('D5').Value = '''211.56'
$ws.Range('E5').Value = '  -1.01%  '
$ws.Range('D6').Value = '''0.493'
$ws.Range('E6').Value = '  -0.75%  '
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('D8').Value = '''46.22'
$ws.Range('D9').Value = '''23.89'
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').Value = '''0.247'
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('D11').Value = '''0.0590'
$ws.Range('E11').Value = '  -1.74%  '
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').Value = '''1.795.73'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('D14').Value = '''1.569.16'
$ws.Range('E14').Value = '  -0.92%  '
$ws.Range('D15').Value = '''0.519'
$ws.Range('E15').Value = '  -2.12%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '''28.454.46'
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = '''3.67'
$ws.Range('E17').Value = '  -2.55%  '
$ws.Range('D18').Value = '''61.86'
$ws.Range('E18').Value = '  -3.25%  '
$ws.Range('D19').Value = '''226.51'
$ws.Range('E19').Value = '  -3.09%  '
$ws.Range('D20').Value = '''7.32'
$ws.Range('E20').Value = '  -2.13%  '
$ws.Range('D21').Value = '0.0₃0691'
$ws.Range('E21').Value = '  -2.47%  '
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').Value = '''3.86'
$ws.Range('E23').Value = '  -6.50%  '
$ws.Range('D24').Value = '''9.10'
$ws.Range('E24').Value = '  -3.00%  '
$ws.Range('E25').Value = '  +7.14%  '
$ws.Range('D26').Value = '''149.90'
$ws.Range('E26').Value = '  -1.08%  '
$ws.Range('D27').Value = '''14.95'
$ws.Range('E27').Value = '  -2.45%  '
$ws.Range('D28').Value = '''6.42'
$ws.Range('E28').Value = '  -2.92%  '
$ws.Range('E29').Value = '  -2.66%  '
$ws.Range('E30').Value = '  +0.53%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '''0.0464'
$ws.Range('E31').Value = '  -1.91%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '''1.11'
$ws.Range('E32').Value = '  -3.60%  '
$ws.Range('E33').Value = '  -1.21%  '
$ws.Range('D34').Value = '''3.13'
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('D35').Value = '''1.394.40'
$ws.Range('E35').Value = '  -1.57%  '
$ws.Range('E36').Value = '  -0.96%  '
$ws.Range('D37').Value = '''1.54'
$ws.Range('E37').Value = '  -3.58%  '
$ws.Range('E38').Value = '  +1.64%  '
$ws.Range('D39').Value = '''2.59'
$ws.Range('E39').Value = '  +1.85%  '
$ws.Range('D40').Value = '''0.0165'
$ws.Range('E40').Value = '  -1.02%  '
$ws.Range('D41').Value = '''0.533'
$ws.Range('E41').Value = '  -1.89%  '
$ws.Range('E42').Value = '  +0.39%  '
$ws.Range('D43').Value = '''0.786'
$ws.Range('E43').Value = '  -3.26%  '
$ws.Range('D44').Value = '''5.57'
$ws.Range('E44').Value = '  -1.91%  '
$ws.Range('D45').Value = '''1.85'
$ws.Range('E45').Value = '  +1.21%  '
$ws.Range('E46').Value = '  +0.80%  '
$ws.Range('D47').Value = '''62.81'
$ws.Range('E47').Value = '  -2.37%  '
$ws.Range('D48').Value = '''1.708.29'
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('D49').Value = '''85.92'
$ws.Range('E49').Value = '  -1.66%  '
$ws.Range('E50').Value = '  -0.50%  '
$ws.Range('D51').Value = '''0.0517'
$ws.Range('E51').Value = '  -1.51%  '
